$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Inflammatory-Mac -> Dkk4 -> Kremen2 -> FAPs (updated TPM values)
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("B2").Value = "Dkk4"
$ws.Range("C2").Value = "Kremen2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.047385
$ws.Range("H2").Value = 0.142155
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03371666666666667
$ws.Range("N2").Value = 0.10115
$ws.Range("O2").Value = 0.5286097276732288
$ws.Range("P2").Value = 0.5286097276732288
$ws.Range("Q2").Value = 0.00159766425
$ws.Range("R2").Value = 0.01437897825
$ws.Range("S2").Value = 0.5286097276732288
$ws.Range("T2").Value = 0.5286097276732288

# Row 3: Inflammatory-Mac -> Dkk4 -> Kremen2 -> MuSCs (updated TPM values)
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Dkk4"
$ws.Range("C3").Value = "Kremen2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.047385
$ws.Range("H3").Value = 0.142155
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.030067
$ws.Range("N3").Value = 0.090201
$ws.Range("O3").Value = 0.4713902723267712
$ws.Range("P3").Value = 0.4713902723267713
$ws.Range("Q3").Value = 0.001424724795
$ws.Range("R3").Value = 0.012822523155
$ws.Range("S3").Value = 0.4713902723267712
$ws.Range("T3").Value = 0.4713902723267713

# Remove the now-obsolete rows 4-7 (old ECs/FAPs/MuSCs/Inflammatory-Mac combinations)
$ws.Range("A4:T7").EntireRow.Delete()
